$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("X-Box Controller")
$ws2 = $wb.Worksheets.Item("Joysticks")

# --- Sheet "X-Box Controller" ---
# New shared strings are introduced in this order: Winch Down, Winch Up,
# Weave Lift Down, Weave Lift Up -- write cells in that order so the
# regenerated shared-string table lines up with the target workbook.

# Row 17: button swaps from "B" to "X", now maps to Winch Down
$ws1.Range("A17").Value = "X"
$ws1.Range("B17").Value = "Winch Down"
$ws1.Range("C17").Value = "Winch Down"

# Row 18: button swaps from "X" to "B", now maps to Winch Up
$ws1.Range("A18").Value = "B"
$ws1.Range("B18").Value = "Winch Up"
$ws1.Range("C18").Value = "Winch Up"

# POV Down (row 13) -> Weave Lift Down
$ws1.Range("B13").Value = "Weave Lift Down"
$ws1.Range("C13").Value = "Weave Lift Down"

# POV Up (row 12) -> Weave Lift Up
$ws1.Range("B12").Value = "Weave Lift Up"
$ws1.Range("C12").Value = "Weave Lift Up"

# --- Sheet "Joysticks" ---
# Button 4 (row 11) -> Winch Down
$ws2.Range("C11").Value = "Winch Down"
$ws2.Range("D11").Value = "Winch Down"

# Button 5 (row 12) -> Winch Up
$ws2.Range("C12").Value = "Winch Up"
$ws2.Range("D12").Value = "Winch Up"

# Button 10 (row 17) -> Weave Lift Down
$ws2.Range("C17").Value = "Weave Lift Down"
$ws2.Range("D17").Value = "Weave Lift Down"

# Button 11 (row 18) -> Weave Lift Up
$ws2.Range("C18").Value = "Weave Lift Up"
$ws2.Range("D18").Value = "Weave Lift Up"

# --- Active sheet / selection changes ---
# Joysticks loses tab selection, its cursor moves to E15
$ws2.Range("E15").Select()

# X-Box Controller becomes the active/selected tab with cursor at D7
$ws1.Activate()
$ws1.Range("D7").Select()
